$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.031.11'
$ws.Cells.Item(2, 5).Value = '  -0.21%  '
$ws.Cells.Item(3, 4).Value = '1.649.58'
$ws.Cells.Item(3, 5).Value = '  -0.69%  '
$ws.Cells.Item(4, 4).Value = '''1.000'
$ws.Cells.Item(4, 5).Value = '  -0.23%  '
$ws.Cells.Item(5, 4).Value = '''214.13'
$ws.Cells.Item(5, 5).Value = '  +3.09%  '
$ws.Cells.Item(6, 4).Value = '''0.5229'
$ws.Cells.Item(6, 5).Value = '  +1.27%  '
$ws.Cells.Item(7, 5).Value = '  -0.20%  '
$ws.Cells.Item(8, 4).Value = '''0.2619'
$ws.Cells.Item(8, 5).Value = '  +1.55%  '
$ws.Cells.Item(9, 4).Value = '''0.06366'
$ws.Cells.Item(9, 5).Value = '  +1.19%  '
$ws.Cells.Item(10, 4).Value = '''20.75'
$ws.Cells.Item(10, 5).Value = '  -0.93%  '
$ws.Cells.Item(11, 4).Value = '''0.07718'
$ws.Cells.Item(11, 5).Value = '  +2.60%  '
$ws.Cells.Item(12, 4).Value = '1.639.86'
$ws.Cells.Item(12, 5).Value = '  -2.92%  '
$ws.Cells.Item(13, 4).Value = '''4.442'
$ws.Cells.Item(13, 5).Value = '  +0.89%  '
$ws.Cells.Item(14, 4).Value = '1.872.80'
$ws.Cells.Item(14, 5).Value = '  -0.58%  '
$ws.Cells.Item(15, 4).Value = '''0.5505'
$ws.Cells.Item(15, 5).Value = '  +2.43%  '
$ws.Cells.Item(16, 4).Value = '0.0₅8266'
$ws.Cells.Item(16, 5).Value = '  +4.37%  '
$ws.Cells.Item(17, 4).Value = '''64.84'
$ws.Cells.Item(17, 5).Value = '  -1.94%  '
$ws.Cells.Item(18, 4).Value = '26.038.00'
$ws.Cells.Item(18, 5).Value = '  -0.18%  '
$ws.Cells.Item(19, 5).Value = '  -0.18%  '
$ws.Cells.Item(20, 4).Value = '''4.736'
$ws.Cells.Item(20, 5).Value = '  +0.92%  '
$ws.Cells.Item(21, 4).Value = '''190.51'
$ws.Cells.Item(21, 5).Value = '  +1.84%  '
$ws.Cells.Item(22, 4).Value = '''10.24'
$ws.Cells.Item(22, 5).Value = '  +0.78%  '
$ws.Cells.Item(23, 4).Value = '''6.327'
$ws.Cells.Item(23, 5).Value = '  +2.52%  '
$ws.Cells.Item(24, 4).Value = '''1.001'
$ws.Cells.Item(24, 5).Value = '  -0.26%  '
$ws.Cells.Item(25, 4).Value = '''143.50'
$ws.Cells.Item(25, 5).Value = '  -3.28%  '
$ws.Cells.Item(26, 4).Value = '''0.1245'
$ws.Cells.Item(26, 5).Value = '  +2.78%  '
$ws.Cells.Item(27, 4).Value = '''7.398'
$ws.Cells.Item(27, 5).Value = '  +0.26%  '
$ws.Cells.Item(28, 4).Value = '''15.98'
$ws.Cells.Item(28, 5).Value = '  +2.68%  '
$ws.Cells.Item(29, 4).Value = '''1.417'
$ws.Cells.Item(29, 5).Value = '  +2.77%  '
$ws.Cells.Item(30, 4).Value = '''0.05940'
$ws.Cells.Item(30, 5).Value = '  -3.44%  '
$ws.Cells.Item(31, 4).Value = '''1.259'
$ws.Cells.Item(31, 5).Value = '  +0.06%  '
$ws.Cells.Item(32, 4).Value = '''3.427'
$ws.Cells.Item(32, 5).Value = '  -1.07%  '
$ws.Cells.Item(33, 4).Value = '''3.409'
$ws.Cells.Item(33, 5).Value = '  +0.36%  '
$ws.Cells.Item(34, 5).Value = '  +1.43%  '
$ws.Cells.Item(35, 4).Value = '''0.9923'
$ws.Cells.Item(35, 5).Value = '  +0.71%  '
$ws.Cells.Item(36, 4).Value = '''2.398'
$ws.Cells.Item(36, 5).Value = '  +0.48%  '
$ws.Cells.Item(37, 4).Value = '''2.759'
$ws.Cells.Item(37, 5).Value = '  +0.45%  '
$ws.Cells.Item(38, 4).Value = '''0.5630'
$ws.Cells.Item(38, 5).Value = '  -3.98%  '
$ws.Cells.Item(39, 5).Value = '  +0.84%  '
$ws.Cells.Item(40, 4).Value = '''5.871'
$ws.Cells.Item(40, 5).Value = '  -1.77%  '
$ws.Cells.Item(41, 4).Value = '''0.8563'
$ws.Cells.Item(41, 5).Value = '  +1.17%  '
$ws.Cells.Item(42, 4).Value = '''1.001'
$ws.Cells.Item(42, 5).Value = '  -0.17%  '
$ws.Cells.Item(43, 4).Value = '1.026.65'
$ws.Cells.Item(43, 5).Value = '  -7.06%  '
$ws.Cells.Item(44, 4).Value = '''99.12'
$ws.Cells.Item(44, 5).Value = '  -0.79%  '
$ws.Cells.Item(45, 4).Value = '1.795.95'
$ws.Cells.Item(45, 5).Value = '  -0.89%  '
$ws.Cells.Item(46, 4).Value = '0.0₈106'
$ws.Cells.Item(46, 5).Value = '  -2.31%  '
$ws.Cells.Item(47, 4).Value = '''55.66'
$ws.Cells.Item(48, 5).Value = '  -0.03%  '
$ws.Cells.Item(49, 4).Value = '''8.033'
$ws.Cells.Item(49, 5).Value = '  +0.07%  '
$ws.Cells.Item(50, 5).Value = '  -1.64%  '

# Row 51: Mantle -> Aptos (coin name, link, price, volume all change)
$ws.Cells.Item(51, 2).Value = 'Aptos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(51, 4).Value = '''5.971'
$ws.Cells.Item(51, 5).Value = '  +2.11%  '
